$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 41 (everything from row 41 down shifts down by one row).
#    This automatically carries merged ranges (B41:E41 -> B42:E42, B44:F44 -> B45:F45,
#    B45:F45 -> B46:F46) and updates formula references that point at the shifted rows
#    (e.g. F42's "=F20+F32+F26+F41" becomes "=F20+F32+F26+F42").
$ws.Rows("41:41").Insert()

# 2. The new row 41 should look like the data-entry rows above it (row 40 pattern):
#    copy formatting from row 40 down into the freshly inserted row 41.
$ws.Range("A40:G40").Copy()
$ws.Range("A41:G41").PasteSpecial(-4122)
$ws.Rows(41).RowHeight = 12.75

# 3. New row 41 is a blank entry row (description blank, merit/demerit 0).
$ws.Range("C41").Value2 = ""
$ws.Range("D41").Value2 = 0
$ws.Range("E41").Value2 = 0

# 4. Row 40's "D" cell now matches the "C" cell styling (reading order alignment).
$ws.Range("C40").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D40").Value2 = 100

# 5. The row that used to be the (blank) totals row is now row 42 and gets the
#    "JUMLAH" label plus an updated total formula that also includes the new row 41.
$ws.Range("B42").Value2 = "JUMLAH"
$ws.Range("F42").Formula = "=D34-E34+D35-E35+D36-E36+D37-E37+D38-E38+D39-E39+D41-E41+D40-E40"

$excel.CutCopyMode = $false
